## "Got 3rd Alt (no regex) working"
## Adds a new "Alt3" worksheet (copied from "Alt2") that extracts the
## bracketed snippet using TEXTSPLIT/MAP/INDEX instead of
## BYROW/REGEXEXTRACT, and makes it the active sheet.

$wb = $excel.ActiveWorkbook

# Keep Alt2 as-is, but first move its own selection/active-cell to where it
# ends up after the new sheet becomes the front-most tab.
$alt2 = $wb.Worksheets.Item("Alt2")
$alt2.Activate() | Out-Null
$alt2.Range("E11").Select() | Out-Null

# Duplicate "Alt2" right after itself - this brings along all of the
# formatting, column widths, styles and cached values of the sheet.
$alt2.Copy($null, $alt2)
$alt3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$alt3.Name = "Alt3"

# Recreate the (hidden, vestigial) AutoFilter defined name that every other
# sheet in this workbook carries, scoped to the new sheet.
$alt3.Names.Add("_xlnm._FilterDatabase", "=Alt3!`$B`$2:`$C`$13", $false) | Out-Null
$filterName = $wb.Names.Item("Alt3!_FilterDatabase")
$filterName.Visible = $false

# Label the new approach.
$alt3.Range("B9").Value = "No Regex"

# Swap the regex-based extraction formula for a TEXTSPLIT/MAP/INDEX based
# one that needs no regex at all.
$formula = '=_xlfn.MAP(B3:B7,_xlfn.LAMBDA(_xlpm.x,INDEX(_xlfn.TEXTSPLIT(_xlpm.x,{"(",")","[","]","{","}"}),,2)))'
$alt3.Range("B11:B15").FormulaArray = $formula

# Make the new sheet the active / front-most tab with its own selection.
$alt3.Activate() | Out-Null
$alt3.Range("F18").Select() | Out-Null
